# Daily Slovak COVID stats update (commit: "Updated: st 23. 09. 2021")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Corrections to existing rows (columns F = AgTests, G = AgPosit) ----
$updates = @(
    @(376, "F", 223078),
    @(400, "F", 150056),
    @(503, "F", 7861),
    @(527, "F", 11623),
    @(530, "F", 12867),
    @(532, "F", 10314),
    @(532, "G", 54),
    @(534, "F", 16797),
    @(535, "F", 10145),
    @(536, "F", 7957),
    @(536, "G", 40),
    @(537, "F", 13997),
    @(537, "G", 52),
    @(538, "F", 11273),
    @(539, "F", 10622),
    @(539, "G", 48),
    @(540, "F", 12472),
    @(541, "F", 16583),
    @(542, "F", 10310),
    @(544, "F", 14350),
    @(545, "F", 16661),
    @(546, "F", 3879),
    @(546, "G", 52),
    @(547, "F", 14048),
    @(548, "F", 17062),
    @(549, "F", 10665),
    @(549, "G", 72),
    @(551, "F", 17738),
    @(551, "G", 198),
    @(552, "G", 175),
    @(553, "F", 15341),
    @(554, "F", 17232),
    @(555, "F", 21492),
    @(557, "F", 10798),
    @(558, "F", 24727),
    @(558, "G", 288),
    @(559, "F", 22341),
    @(559, "G", 266),
    @(560, "F", 5930),
    @(560, "G", 90),
    @(561, "F", 23292),
    @(561, "G", 374),
    @(562, "F", 26081),
    @(562, "G", 261),
    @(563, "F", 13495),
    @(563, "G", 161),
    @(564, "F", 13621),
    @(564, "G", 190)
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $val = $u[2]
    $ws.Range("$col$row").Value = $val
}

# ---- Append three new daily rows (565-567) ----
$newRows = @(
    @(565, 44459, 403802, 10475, 880, 12580, 27343, 353),
    @(566, 44460, 404982, 10898, 1180, 12589, 24616, 325),
    @(567, 44461, 405931, 9925, 949, 12592, 15392, 422)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Range("A$rowNum").Value = $r[1]
    $ws.Range("B$rowNum").Value = $r[2]
    $ws.Range("C$rowNum").Value = $r[3]
    $ws.Range("D$rowNum").Value = $r[4]
    $ws.Range("E$rowNum").Value = $r[5]
    $ws.Range("F$rowNum").Value = $r[6]
    $ws.Range("G$rowNum").Value = $r[7]
}
